$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-5 (existing MuSCs block becomes ECs block with new TPM values)
# Add new rows 6-9 (new MuSCs block with new TPM values)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il1rapl1"
$ws.Range("C2").Value = "Ptprs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.03825666666666667
$ws.Range("H2").Value = 0.11477
$ws.Range("I2").Value = 0.4331352014340976
$ws.Range("J2").Value = 0.4331352014340976
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.556762333333333
$ws.Range("N2").Value = 10.670287
$ws.Range("O2").Value = 0.04280930450251701
$ws.Range("P2").Value = 0.04280930450251701
$ws.Range("Q2").Value = 0.1360698709988889
$ws.Range("R2").Value = 1.22462883899
$ws.Range("S2").Value = 0.01854221672895133
$ws.Range("T2").Value = 0.01854221672895133

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il1rapl1"
$ws.Range("C3").Value = "Ptprs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.03825666666666667
$ws.Range("H3").Value = 0.11477
$ws.Range("I3").Value = 0.4331352014340976
$ws.Range("J3").Value = 0.4331352014340976
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 47.24901333333333
$ws.Range("N3").Value = 141.74704
$ws.Range("O3").Value = 0.5686906263805706
$ws.Range("P3").Value = 0.5686906263805704
$ws.Range("Q3").Value = 1.807589753422222
$ws.Range("R3").Value = 16.2683077808
$ws.Range("S3").Value = 0.2463199290110316
$ws.Range("T3").Value = 0.2463199290110315

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il1rapl1"
$ws.Range("C4").Value = "Ptprs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.03825666666666667
$ws.Range("H4").Value = 0.11477
$ws.Range("I4").Value = 0.4331352014340976
$ws.Range("J4").Value = 0.4331352014340976
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 24.53173066666666
$ws.Range("N4").Value = 73.595192
$ws.Range("O4").Value = 0.2952646900921413
$ws.Range("P4").Value = 0.2952646900921412
$ws.Range("Q4").Value = 0.9385022428711111
$ws.Range("R4").Value = 8.446520185839999
$ws.Range("S4").Value = 0.127889531019436
$ws.Range("T4").Value = 0.127889531019436

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Il1rapl1"
$ws.Range("C5").Value = "Ptprs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.03825666666666667
$ws.Range("H5").Value = 0.11477
$ws.Range("I5").Value = 0.4331352014340976
$ws.Range("J5").Value = 0.4331352014340976
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.746355333333334
$ws.Range("N5").Value = 23.239066
$ws.Range("O5").Value = 0.09323537902477132
$ws.Range("P5").Value = 0.0932353790247713
$ws.Range("Q5").Value = 0.2963497338688889
$ws.Range("R5").Value = 2.66714760482
$ws.Range("S5").Value = 0.04038352467467876
$ws.Range("T5").Value = 0.04038352467467876

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Il1rapl1"
$ws.Range("C6").Value = "Ptprs"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.05006833333333333
$ws.Range("H6").Value = 0.150205
$ws.Range("I6").Value = 0.5668647985659024
$ws.Range("J6").Value = 0.5668647985659024
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.556762333333333
$ws.Range("N6").Value = 10.670287
$ws.Range("O6").Value = 0.04280930450251701
$ws.Range("P6").Value = 0.04280930450251701
$ws.Range("Q6").Value = 0.1780811620927777
$ws.Range("R6").Value = 1.602730458835
$ws.Range("S6").Value = 0.02426708777356569
$ws.Range("T6").Value = 0.02426708777356568

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Il1rapl1"
$ws.Range("C7").Value = "Ptprs"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.05006833333333333
$ws.Range("H7").Value = 0.150205
$ws.Range("I7").Value = 0.5668647985659024
$ws.Range("J7").Value = 0.5668647985659024
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 47.24901333333333
$ws.Range("N7").Value = 141.74704
$ws.Range("O7").Value = 0.5686906263805706
$ws.Range("P7").Value = 0.5686906263805704
$ws.Range("Q7").Value = 2.365679349244444
$ws.Range("R7").Value = 21.2911141432
$ws.Range("S7").Value = 0.322370697369539
$ws.Range("T7").Value = 0.3223706973695389

# Row 8
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Il1rapl1"
$ws.Range("C8").Value = "Ptprs"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.05006833333333333
$ws.Range("H8").Value = 0.150205
$ws.Range("I8").Value = 0.5668647985659024
$ws.Range("J8").Value = 0.5668647985659024
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 24.53173066666666
$ws.Range("N8").Value = 73.595192
$ws.Range("O8").Value = 0.2952646900921413
$ws.Range("P8").Value = 0.2952646900921412
$ws.Range("Q8").Value = 1.228262868262222
$ws.Range("R8").Value = 11.05436581436
$ws.Range("S8").Value = 0.1673751590727053
$ws.Range("T8").Value = 0.1673751590727053

# Row 9
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Il1rapl1"
$ws.Range("C9").Value = "Ptprs"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.05006833333333333
$ws.Range("H9").Value = 0.150205
$ws.Range("I9").Value = 0.5668647985659024
$ws.Range("J9").Value = 0.5668647985659024
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 7.746355333333334
$ws.Range("N9").Value = 23.239066
$ws.Range("O9").Value = 0.09323537902477132
$ws.Range("P9").Value = 0.0932353790247713
$ws.Range("Q9").Value = 0.3878471009477777
$ws.Range("R9").Value = 3.490623908529999
$ws.Range("S9").Value = 0.05285185435009256
$ws.Range("T9").Value = 0.05285185435009255

